# Fill in Time Log entries for rows 51 and 52 on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 51: Date, Start Time, Stop Time, Interruption (minutes), Activity
# (set the Interruption value first so the dependent formula in column E
#  recalculates against the correct set of inputs)
$ws.Cells.Item(51, 4).Value = 20
$ws.Cells.Item(51, 1).Value = 41899
$ws.Cells.Item(51, 2).Value = 0.95138888888888884
$ws.Cells.Item(51, 3).Value = 1.0506944444444444
$ws.Cells.Item(51, 6).Value = "Testing"

# Row 52: Date, Start Time, Stop Time, Interruption (minutes), Activity
$ws.Cells.Item(52, 4).Value = 20
$ws.Cells.Item(52, 1).Value = 41900
$ws.Cells.Item(52, 2).Value = 0.69930555555555562
$ws.Cells.Item(52, 3).Value = 0.78333333333333333
$ws.Cells.Item(52, 6).Value = "Testing"

# Move the active selection to where the user ended up after entering the
# two new rows of data (C53)
$ws.Range("C53").Select()

$wb.Save()
